{"js": "// The document currently spells out the id tag in three separate runs:\n//   <id>  (Courier New, color 7f6000, sz 18)  +  p161r_1  (plain)  +  </id>  (Courier New, color 7f6000, sz 18)\n// The commit collapses these three runs into a single run containing the\n// full text \"<id>p161r_1</id>\", keeping the first run's (Courier New /\n// 7f6000 / 18pt) character formatting.\nconst body = context.document.body;\n\n// Locate the run span by its concatenated plain text (search reads across\n// run boundaries within a paragraph).\nconst results = body.search(\"<id>p161r_1</id>\", { matchCase: true, matchWildcards: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find '<id>p161r_1</id>' in the document body.\");\n}\n\n// Replacing the whole matched range's text collapses it into a single run\n// that carries the formatting of the range's (first) run - exactly the\n// <id> run's Courier New / 7f6000 / 18pt properties - so the 2nd/3rd runs\n// disappear and the 1st run's text becomes the full \"<id>p161r_1</id>\".\nconst target = results.items[0];\ntarget.insertText(\"<id>p161r_1</id>\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The \"<id>p161r_1</id>\" text is currently split across three runs:\n#   <id>       - Courier New, color 7f6000, sz 18/szCs 18  (tag-markup style)\n#   p161r_1    - default/plain style\n#   </id>      - Courier New, color 7f6000, sz 18/szCs 18  (tag-markup style)\n# The edit merges them into ONE run containing the full text\n# \"<id>p161r_1</id>\", carrying the \"<id>\" run's (Courier New / 7f6000 / 18pt)\n# character formatting - i.e. the 2nd/3rd runs disappear and their text gets\n# folded into the first run.\n\n$d = $word.ActiveDocument\n\n# Remove the \"p161r_1</id>\" tail (the 2nd + 3rd runs) ...\n$tail = $d.Content\n$foundTail = $tail.Find.Execute(\"p161r_1</id>\")\nif (-not $foundTail) {\n    throw \"Could not find 'p161r_1</id>' in the document.\"\n}\n$tail.Delete()\n\n# ... then re-append that same text right after the remaining \"<id>\" run,\n# which adopts the \"<id>\" run's formatting so everything collapses back into\n# a single run reading \"<id>p161r_1</id>\".\n$head = $d.Content\n$foundHead = $head.Find.Execute(\"<id>\")\nif (-not $foundHead) {\n    throw \"Could not find '<id>' in the document.\"\n}\n$head.InsertAfter(\"p161r_1</id>\")\n"}
